$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# Row 6
$wsALC.Range("H6").Value = 86.8
$wsALC.Range("I6").Value = 83.5
$wsALC.Range("J6").Value = 100
$wsALC.Range("K6").Value = 250.5
$wsALC.Range("L6").Value = 300
$wsALC.Range("M6").Value = -138.5
$wsALC.Range("N6").Value = -524
# Row 12
$wsALC.Range("H12").Value = 245
$wsALC.Range("I12").Value = 250
$wsALC.Range("J12").Value = 240
$wsALC.Range("K12").Value = 250
$wsALC.Range("L12").Value = 240
$wsALC.Range("M12").Value = -80
$wsALC.Range("N12").Value = -580
# Row 29
$wsALC.Range("H29").Value = 7474.75
$wsALC.Range("J29").Value = 7466.6665
$wsALC.Range("L29").Value = 22399.9995
$wsALC.Range("N29").Value = -22961.9995
# Row 38
$wsALC.Range("H38").Value = 3709.9
$wsALC.Range("I38").Value = 176.5
$wsALC.Range("J38").Value = 6065.5
$wsALC.Range("K38").Value = 529.5
$wsALC.Range("L38").Value = 18196.5
$wsALC.Range("M38").Value = -157.5
$wsALC.Range("N38").Value = -18940.5
# Row 62
$wsALC.Range("H62").Value = 2701.1428
$wsALC.Range("I62").Value = 2680.5789
$wsALC.Range("K62").Value = 2680.5789
$wsALC.Range("M62").Value = -2056.5789
# Row 65
$wsALC.Range("H65").Value = 2701.1428
$wsALC.Range("I65").Value = 2680.5789
$wsALC.Range("K65").Value = 13402.8945
$wsALC.Range("M65").Value = -10282.8945
# Row 87
$wsALC.Range("H87").Value = 200000
$wsALC.Range("J87").Value = 200000
$wsALC.Range("L87").Value = 200000
$wsALC.Range("N87").Value = -202496
# Row 90
$wsALC.Range("H90").Value = 200000
$wsALC.Range("J90").Value = 200000
$wsALC.Range("L90").Value = 600000
$wsALC.Range("N90").Value = -612480
# Row 99
$wsALC.Range("H99").Value = 291
$wsALC.Range("I99").Value = 303.81818
$wsALC.Range("J99").Value = 150
$wsALC.Range("K99").Value = 911.45454
$wsALC.Range("L99").Value = 450
$wsALC.Range("M99").Value = 586.54546
$wsALC.Range("N99").Value = -3446
# Row 132
$wsALC.Range("H132").Value = 6634.0713
$wsALC.Range("I132").Value = 6857.5557
$wsALC.Range("K132").Value = 20572.6671
$wsALC.Range("M132").Value = -18042.6671
# Row 137
$wsALC.Range("H137").Value = 8253.799999999999
$wsALC.Range("I137").Value = 2057.5
$wsALC.Range("J137").Value = 15335.286
$wsALC.Range("K137").Value = 6172.5
$wsALC.Range("L137").Value = 46005.858
$wsALC.Range("M137").Value = -3622.5
$wsALC.Range("N137").Value = -51105.858

# --- ARM ---
# Row 32
$wsARM.Range("H32").Value = 3491.8364
$wsARM.Range("I32").Value = 3491.8364
$wsARM.Range("K32").Value = 3491.8364
$wsARM.Range("M32").Value = -3204.8364
# Row 45
$wsARM.Range("H45").Value = 18430.133
$wsARM.Range("I45").Value = 24251.38
$wsARM.Range("K45").Value = 24251.38
$wsARM.Range("M45").Value = -23874.38

# --- BSM ---
# Row 16
$wsBSM.Range("H16").Value = 0
$wsBSM.Range("J16").Value = 0
$wsBSM.Range("L16").Value = 0
$wsBSM.Range("N16").Value = $null

# --- CRP ---
# Row 14
$wsCRP.Range("H14").Value = 0
$wsCRP.Range("I14").Value = 0
$wsCRP.Range("K14").Value = 0
$wsCRP.Range("M14").Value = $null
# Row 22
$wsCRP.Range("H22").Value = 1655.875
$wsCRP.Range("I22").Value = 1229.6
$wsCRP.Range("K22").Value = 1229.6
$wsCRP.Range("M22").Value = -879.5999999999999
# Row 31
$wsCRP.Range("H31").Value = 3586.8438
$wsCRP.Range("J31").Value = 5056.25
$wsCRP.Range("L31").Value = 5056.25
$wsCRP.Range("N31").Value = -5646.25
# Row 34
$wsCRP.Range("H34").Value = 3586.8438
$wsCRP.Range("J34").Value = 5056.25
$wsCRP.Range("L34").Value = 5056.25
$wsCRP.Range("N34").Value = -5460.25
# Row 86
$wsCRP.Range("H86").Value = 1109401.4
$wsCRP.Range("I86").Value = 3097.2144
$wsCRP.Range("J86").Value = 3322009.8
$wsCRP.Range("K86").Value = 3097.2144
$wsCRP.Range("L86").Value = 3322009.8
$wsCRP.Range("M86").Value = -1974.2144
$wsCRP.Range("N86").Value = -3324255.8
# Row 89
$wsCRP.Range("H89").Value = 1109401.4
$wsCRP.Range("I89").Value = 3097.2144
$wsCRP.Range("J89").Value = 3322009.8
$wsCRP.Range("K89").Value = 15486.072
$wsCRP.Range("L89").Value = 16610049
$wsCRP.Range("M89").Value = -9870.072
$wsCRP.Range("N89").Value = -16621281
# Row 122
$wsCRP.Range("H122").Value = 4497.25
$wsCRP.Range("I122").Value = 2982.8572
$wsCRP.Range("K122").Value = 8948.571599999999
$wsCRP.Range("M122").Value = -6498.571599999999
# Row 132
$wsCRP.Range("H132").Value = 2042.9318
$wsCRP.Range("I132").Value = 1494.7273
$wsCRP.Range("K132").Value = 4484.1819
$wsCRP.Range("M132").Value = -1954.1819

# --- CUL ---
# Row 5
$wsCUL.Range("H5").Value = 2008.95
$wsCUL.Range("I5").Value = 577.3333
$wsCUL.Range("J5").Value = 2261.5881
$wsCUL.Range("K5").Value = 1731.9999
$wsCUL.Range("L5").Value = 6784.7643
$wsCUL.Range("M5").Value = -1619.9999
$wsCUL.Range("N5").Value = -7008.7643
# Row 56
$wsCUL.Range("H56").Value = 6315.8335
$wsCUL.Range("I56").Value = 6315.8335
$wsCUL.Range("K56").Value = 6315.8335
$wsCUL.Range("M56").Value = -5785.8335
# Row 107
$wsCUL.Range("H107").Value = 491.6
$wsCUL.Range("I107").Value = 299
$wsCUL.Range("J107").Value = 620
$wsCUL.Range("K107").Value = 897
$wsCUL.Range("L107").Value = 1860
$wsCUL.Range("M107").Value = 1023
$wsCUL.Range("N107").Value = -5700
# Row 109
$wsCUL.Range("H109").Value = 2215.3333
$wsCUL.Range("I109").Value = 2215.3333
$wsCUL.Range("K109").Value = 6645.999899999999
$wsCUL.Range("M109").Value = -5605.999899999999
# Row 112
$wsCUL.Range("H112").Value = 6249.5
$wsCUL.Range("I112").Value = 4999
$wsCUL.Range("J112").Value = 7500
$wsCUL.Range("K112").Value = 14997
$wsCUL.Range("L112").Value = 22500
$wsCUL.Range("M112").Value = -13889
$wsCUL.Range("N112").Value = -24716
# Row 115
$wsCUL.Range("H115").Value = 6876.1113
$wsCUL.Range("J115").Value = 8285.286
$wsCUL.Range("L115").Value = 24855.858
$wsCUL.Range("N115").Value = -27205.858
# Row 126
$wsCUL.Range("H126").Value = 0
$wsCUL.Range("I126").Value = 0
$wsCUL.Range("K126").Value = 0
$wsCUL.Range("M126").Value = $null
# Row 130
$wsCUL.Range("H130").Value = 11510
$wsCUL.Range("J130").Value = 7250
$wsCUL.Range("L130").Value = 21750
$wsCUL.Range("N130").Value = -31790
# Row 131
$wsCUL.Range("H131").Value = 7694.4165
$wsCUL.Range("I131").Value = 11482
$wsCUL.Range("J131").Value = 2391.8
$wsCUL.Range("K131").Value = 34446
$wsCUL.Range("L131").Value = 7175.400000000001
$wsCUL.Range("M131").Value = -29406
$wsCUL.Range("N131").Value = -17255.4
# Row 132
$wsCUL.Range("H132").Value = 2234.5881
$wsCUL.Range("I132").Value = 1329.25
$wsCUL.Range("K132").Value = 11963.25
$wsCUL.Range("M132").Value = -9433.25
# Row 135
$wsCUL.Range("H135").Value = 2008.95
$wsCUL.Range("I135").Value = 577.3333
$wsCUL.Range("J135").Value = 2261.5881
$wsCUL.Range("K135").Value = 5195.9997
$wsCUL.Range("L135").Value = 20354.2929
$wsCUL.Range("M135").Value = -2660.9997
$wsCUL.Range("N135").Value = -25424.2929
# Row 141
$wsCUL.Range("H141").Value = 11447.333
$wsCUL.Range("I141").Value = 11447.333
$wsCUL.Range("K141").Value = 34341.999
$wsCUL.Range("M141").Value = -29161.999

# --- LTW ---
# Row 114
$wsLTW.Range("H114").Value = 119999.5
$wsLTW.Range("J114").Value = 119999.5
$wsLTW.Range("L114").Value = 119999.5
$wsLTW.Range("N114").Value = -128677.5
# Row 122
$wsLTW.Range("H122").Value = 4499.7
$wsLTW.Range("I122").Value = 4724.75
$wsLTW.Range("J122").Value = 3599.5
$wsLTW.Range("K122").Value = 14174.25
$wsLTW.Range("L122").Value = 10798.5
$wsLTW.Range("M122").Value = -11724.25
$wsLTW.Range("N122").Value = -15698.5
# Row 132
$wsLTW.Range("H132").Value = 3192.0967
$wsLTW.Range("I132").Value = 2453.283
$wsLTW.Range("K132").Value = 7359.849
$wsLTW.Range("M132").Value = -4829.849

# --- WVR ---
# Row 23
$wsWVR.Range("H23").Value = 12999.5
$wsWVR.Range("I23").Value = 10999
$wsWVR.Range("K23").Value = 10999
$wsWVR.Range("M23").Value = -10770
# Row 107
$wsWVR.Range("H107").Value = 428.07693
$wsWVR.Range("I107").Value = 435.83334
$wsWVR.Range("J107").Value = 335
$wsWVR.Range("K107").Value = 1307.50002
$wsWVR.Range("L107").Value = 1005
$wsWVR.Range("M107").Value = 612.4999800000001
$wsWVR.Range("N107").Value = -4845
# Row 136
$wsWVR.Range("H136").Value = 23811712
$wsWVR.Range("I136").Value = 27027964
$wsWVR.Range("K136").Value = 81083892
$wsWVR.Range("M136").Value = -81081342
